$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E (publication_date) and M (cited_by_count) hold numeric- or
# date-looking text but must keep their "Text" cell type, matching the
# source file (every cell there is a literal string, t="inlineStr").
# Forcing NumberFormat to "@" before writing prevents Excel/COM from
# silently re-typing "79" as a number or "2023-09-01" as a date serial.

# --- Row 2 ---
$ws.Range("A2").Value = "Mina K. Chung, Kristen K. Patton, Chu‐Pak Lau, Alexander Romeno Janner Dal Forno, Sana M. Al‐Khatib, Vanita Arora, Ulrika Birgersdotter‐Green, Yong‐Mei Cha, Eugene H. Chung, Edmond M. Cronin, Anne B. Curtis, Iwona Cygankiewicz, Gopi Dandamudi, Anne M. Dubin, Douglas P. Ensch, Taya V. Glotzer, Michael R. Gold, Zachary D. Goldberger, Rakesh Gopinathannair, Eiran Z. Gorodeski, Alejandra Gutiérrez, Juan C. Guzmán, Weijian Huang, Peter B. Imrey, Julia H. Indik, Saima Karim, Peter P. Karpawich, Yaariv Khaykin, Erich L. Kiehl, Jordana Kron, Valentina Kutyifa, Mark S. Link, Joseph E. Marine, Wilfried Mullens, Seung Jung Park, Ratika Parkash, Manuel Patete, Rajeev K. Pathak, Carlos Perona, John Rickard, Mark H. Schoenfeld, Swee‐Chong Seow, Win‐Kuang Shen, Morio Shoda, Jagmeet P. Singh, David J. Slotwiner, Arun Raghav Mahankali Sridhar, Uma N Srivatsa, Eric C. Stecker, Tanyanan Tanawuttiwat, W.H. Wilson Tang, Carlos Andres Tapias, Cynthia M. Tracy, Gaurav A. Upadhyay, Niraj Varma, Kevin Vernooy, Pugazhendhi Vijayaraman, Sarah Ann Worsnick, Wojciech Zaręba, Emily P. Zeitler"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "79"

# --- Row 5 ---
$ws.Range("A5").Value = "Mina K. Chung, Kristen K. Patton, Chu‐Pak Lau, Alexander Romeno Janner Dal Forno, Sana M. Al‐Khatib, Vanita Arora, Ulrika Birgersdotter‐Green, Yong‐Mei Cha, Eugene H. Chung, Edmond M. Cronin, Anne B. Curtis, Iwona Cygankiewicz, Gopi Dandamudi, Anne M. Dubin, Douglas P. Ensch, Taya V. Glotzer, Michael R. Gold, Zachary D. Goldberger, Rakesh Gopinathannair, Eiran Z. Gorodeski, Alejandra Gutiérrez, Juan C. Guzmán, Weijian Huang, Peter B. Imrey, Julia H. Indik, Saima Karim, Peter P. Karpawich, Yaariv Khaykin, Erich L. Kiehl, Jordana Kron, Valentina Kutyifa, Mark S. Link, Joseph E. Marine, Wilfried Mullens, Seung‐Jung Park, Ratika Parkash, Manuel Patete, Rajeev K. Pathak, Carlos Perona, John Rickard, Mark H. Schoenfeld, Swee‐Chong Seow, Win‐Kuang Shen, Morio Shoda, Jagmeet P. Singh, David J. Slotwiner, Arun Raghav Mahankali Sridhar, Uma N Srivatsa, Eric C. Stecker, Tanyanan Tanawuttiwat, W.H. Wilson Tang, Carlos Andres Tapias, Cynthia M. Tracy, Gaurav A. Upadhyay, Niraj Varma, Kevin Vernooy, Pugazhendhi Vijayaraman, Sarah Ann Worsnick, Wojciech Zaręba, Emily P. Zeitler"

# --- Row 7 ---
$ws.Range("B7").ClearContents()
$ws.Range("A7").Value = "Julia H. Indik"
$ws.Range("C7").Value = "https://openalex.org/W4386226754"
$ws.Range("D7").Value = "Introducing the 2023 HRS/APHRS/LAHRS guideline on cardiac physiologic pacing for the avoidance and mitigation of heart failure: Are we entering a new age in pacing?"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2023-09-01"
$ws.Range("F7").Value = "Heart Rhythm O2"
$ws.Range("H7").Value = "https://doi.org/10.1016/j.hroo.2023.08.002"
$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = "1"
$ws.Range("O7").Value = "https://pubmed.ncbi.nlm.nih.gov/37744941"
$ws.Range("P7").Value = "https://doi.org/10.1016/j.hroo.2023.08.002"

# --- Row 8 ---
$ws.Range("A8").Value = "Julia H. Indik, Hugh Calkins"
$ws.Range("B8").Value = "University of Arizona College of Medicine, Tucson, Arizona; Johns Hopkins Hospital, Baltimore, Maryland"
$ws.Range("C8").Value = "https://openalex.org/W4321787357"
$ws.Range("D8").Value = "Frank Marcus (March 23, 1928–December 21, 2022)"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2023-03-01"
$ws.Range("F8").Value = "Heart Rhythm"
$ws.Range("G8").Value = "Elsevier BV"
$ws.Range("H8").Value = "https://doi.org/10.1016/j.hrthm.2023.01.015"
$ws.Range("O8").Value = "https://pubmed.ncbi.nlm.nih.gov/36842791"
$ws.Range("P8").Value = "https://doi.org/10.1016/j.hrthm.2023.01.015"

# --- Row 9 ---
$ws.Range("C9").Value = "https://openalex.org/W4366351470"
$ws.Range("D9").Value = "Rhythm Control Treatment for Atrial Fibrillation Is Not Just for the Healthy"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2023-05-01"
$ws.Range("F9").Value = "Circulation: Arrhythmia and Electrophysiology"
$ws.Range("G9").Value = "Lippincott Williams & Wilkins"
$ws.Range("H9").Value = "https://doi.org/10.1161/circep.123.011949"
$ws.Range("O9").Value = "https://pubmed.ncbi.nlm.nih.gov/37073732"
$ws.Range("P9").Value = "https://doi.org/10.1161/circep.123.011949"
$ws.Range("B9").Value = "University of Arizona College of Medicine, Tucson, AZ."
